$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155; this shifts the existing rows
# 155-180 down to 156-181 (matching the dimension growing from
# A1:T180 to A1:T181).
$ws.Rows(155).Insert()

# Populate the newly inserted row 155 with the new record.
$ws.Range("A155").Value = 11
$ws.Range("B155").Value = "Vega Monumental Concepción"
$ws.Range("C155").Value = "Bíobío"
$ws.Range("D155").Value = 45127
$ws.Range("E155").Value = 8
$ws.Range("F155").Value = "Fruta"
$ws.Range("G155").Value = 100108
$ws.Range("H155").Value = "Tropicales y subtropicales"
$ws.Range("I155").Value = 100108002
$ws.Range("J155").Value = "Mango"
$ws.Range("K155").Value = "Sin especificar"
$ws.Range("L155").Value = "Primera"
$ws.Range("M155").Value = 200
$ws.Range("N155").Value = 8500
$ws.Range("O155").Value = 9000
$ws.Range("P155").Value = 8750
$ws.Range("Q155").Value = "$/bandeja 4 kilos"
$ws.Range("R155").Value = "Brasil"
$ws.Range("S155").Value = 2188
$ws.Range("T155").Value = 4
